# Sync attendance_reports: reorder "Recorded By" (column G) name lists.
# The underlying data source evidently re-ordered how recorder names are
# concatenated; this script re-applies the same literal text substitutions
# to every cell in column G that still holds the old ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Exact (old -> new) text replacements observed for the "Recorded By" column.
$replacements = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
